$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the actual spending figure for "Sosiale tilstellinger" (E48) to
# reflect the new equipment purchase; downstream SUM/diff formulas
# (F48, E78, F78, E80, F80, E82, F82) recalculate automatically.
$ws.Range("E48").Value = 5044.29

# Move the view/selection to where the editor left off.
$ws.Range("G52").Select() | Out-Null
